$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the added columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copy the header style (s=3) from K1 into the new header cells
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 2-6: E and F columns were scaled from fraction (0-1) to percentage number (0-100)
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 62.17228464419475

$ws.Range("E3").Value = 0.3405994550408719
$ws.Range("F3").Value = 0

$ws.Range("E4").Value = 99.65940054495913
$ws.Range("F4").Value = 94.53178400546821

$ws.Range("E5").Value = 91.66666666666666
$ws.Range("F5").Value = 21.85007974481659

$ws.Range("E6").Value = 8.333333333333332
$ws.Range("F6").Value = 26.31578947368421

# New columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes) for rows 2-6
$ws.Range("L2").Value = 91.30338044842225
$ws.Range("M2").Value = 263553
$ws.Range("N2").Value = 317.533734939759

$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0

$ws.Range("L4").Value = 90.1669167946294
$ws.Range("M4").Value = 203646
$ws.Range("N4").Value = 147.2494577006508

$ws.Range("L5").Value = 19.17132323902399
$ws.Range("M5").Value = 2063
$ws.Range("N5").Value = 15.05839416058394

$ws.Range("L6").Value = 25.0794338805401
$ws.Range("M6").Value = 145
$ws.Range("N6").Value = 9.666666666666666
